# Updates the cryptocurrency price/volume table to the latest scrape.
# Each row corresponds to one coin; only the Price (D) and Volume(1h) (E)
# columns change value, except row 51 which was fully replaced with a new coin.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.529.13"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "2.642.74"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'603.64"
$ws.Range("E5").Value = "  +2.27%  "
$ws.Range("D6").Value = "'146.27"
$ws.Range("E6").Value = "  +1.95%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("E9").Value = "  +1.76%  "
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("D11").Value = "'0.370"
$ws.Range("E11").Value = "  +4.97%  "
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").Value = "'27.57"
$ws.Range("E13").Value = "  +0.89%  "
$ws.Range("D14").Value = "3.117.21"
$ws.Range("E14").Value = "  -0.15%  "
$ws.Range("D15").Value = "63.338.38"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").Value = "'0.0000147"
$ws.Range("E16").Value = "  +1.66%  "
$ws.Range("D17").Value = "2.669.57"
$ws.Range("E17").Value = "  +1.69%  "
$ws.Range("D18").Value = "'11.49"
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("D19").Value = "'4.58"
$ws.Range("E19").Value = "  +5.13%  "
$ws.Range("D20").Value = "'344.05"
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("E21").Value = "  +2.93%  "
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("E23").Value = "  -3.33%  "
$ws.Range("D24").Value = "'66.57"
$ws.Range("E24").Value = "  -0.92%  "
$ws.Range("D25").Value = "'1.69"
$ws.Range("E25").Value = "  +1.95%  "
$ws.Range("D26").Value = "'9.13"
$ws.Range("E26").Value = "  +8.55%  "
$ws.Range("D27").Value = "'579.48"
$ws.Range("E27").Value = "  +6.82%  "
$ws.Range("D28").Value = "'1.55"
$ws.Range("E29").Value = "  -1.39%  "
$ws.Range("E30").Value = "  -0.58%  "
$ws.Range("D31").Value = "'7.94"
$ws.Range("E31").Value = "  +2.36%  "
$ws.Range("D32").Value = "'2.06"
$ws.Range("E32").Value = "  +4.22%  "
$ws.Range("E33").Value = "  -2.40%  "
$ws.Range("D34").Value = "0.0₃0826"
$ws.Range("E34").Value = "  +2.70%  "
$ws.Range("D35").Value = "'5.22"
$ws.Range("E35").Value = "  +7.73%  "
$ws.Range("D36").Value = "'166.90"
$ws.Range("E36").Value = "  -4.83%  "
$ws.Range("E37").Value = "  +1.13%  "
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("D39").Value = "'1.95"
$ws.Range("E39").Value = "  +8.70%  "
$ws.Range("D40").Value = "'19.12"
$ws.Range("E40").Value = "  +0.41%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").Value = "'168.94"
$ws.Range("E42").Value = "  -1.86%  "
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("E45").Value = "  +1.15%  "
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").Value = "'0.0246"
$ws.Range("E47").Value = "  +3.32%  "
$ws.Range("D48").Value = "'0.0962"
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("E49").Value = "  +12.18%  "
$ws.Range("D50").Value = "'18.80"
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("B51").Value = "TheGraph"
$ws.Range("C51").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D51").Value = "'0.177"
$ws.Range("E51").Value = "  +1.20%  "
